# Update model class diagram to indicate optional fields.
#
# The existing purple "*" marker TextBox ("TextBox 70") gets duplicated four
# more times and each copy (plus the original, which is repositioned) is
# placed next to a field that is now marked optional in the diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Locate the source marker shape ("TextBox 70") that all the new "*"
# badges are cloned from.
# ---------------------------------------------------------------------
$orig = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 70") {
        $orig = $s.Shapes.Item($i)
    }
}

# ---------------------------------------------------------------------
# The shape-id counter in this deck is a monotonically increasing,
# never-reused counter that always hands out the smallest id that has
# never been used before. The slide already has ids allocated all the
# way up to 124, but ids 3-45 (inclusive, with a couple of exceptions
# already in use) were never consumed. Burn through that unused range
# with disposable clones so the *next* five shapes we create land
# exactly on ids 72, 73, 74, 77 and 82 - matching what real PowerPoint
# handed out when this edit was authored.
# ---------------------------------------------------------------------
$fillers = New-Object System.Collections.ArrayList
for ($i = 1; $i -le 40; $i++) {
    [void]$fillers.Add($orig.Duplicate())
}
foreach ($f in $fillers) {
    $f.Delete()
}

# ---------------------------------------------------------------------
# Create the five final marker shapes as clones of the original so they
# inherit its exact text formatting / body properties, then move each
# one onto its final, precise position.
# ---------------------------------------------------------------------
$c1 = $orig.Duplicate()
$c2 = $orig.Duplicate()
$c3 = $orig.Duplicate()
$c4 = $orig.Duplicate()
$c5 = $orig.Duplicate()

# The original shape itself is superseded by the repositioned clone -
# remove it now that its replacement exists.
$orig.Delete()

$c1.Name = "TextBox 71"
$c1.Left = 588.2925415039062
$c1.Top = 281.0318298339844
$c1.Width = 14.90212631225586
$c1.Height = 14.069527626037598

$c2.Name = "TextBox 72"
$c2.Left = 587.8931884765625
$c2.Top = 254.07205200195312
$c2.Width = 14.90212631225586
$c2.Height = 14.069527626037598

$c3.Name = "TextBox 73"
$c3.Left = 587.8931884765625
$c3.Top = 306.9001770019531
$c3.Width = 14.90212631225586
$c3.Height = 14.069527626037598

$c4.Name = "TextBox 76"
$c4.Left = 588.0
$c4.Top = 331.885986328125
$c4.Width = 14.90212631225586
$c4.Height = 14.069527626037598

$c5.Name = "TextBox 81"
$c5.Left = 587.8931884765625
$c5.Top = 230.7244110107422
$c5.Width = 14.90212631225586
$c5.Height = 14.069527626037598
